$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F (reuse the same formatting as the other header cells)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:50:32.587557",
    "2021-10-05 10:50:32.587569",
    "2021-10-05 10:50:32.587573",
    "2021-10-05 10:50:32.587576",
    "2021-10-05 10:50:32.587580",
    "2021-10-05 10:50:32.587583",
    "2021-10-05 10:50:32.587586",
    "2021-10-05 10:50:32.587589",
    "2021-10-05 10:50:32.587593",
    "2021-10-05 10:50:32.587596",
    "2021-10-05 10:50:32.587599",
    "2021-10-05 10:50:32.587602",
    "2021-10-05 10:50:32.587605",
    "2021-10-05 10:50:32.587608",
    "2021-10-05 10:50:32.587611",
    "2021-10-05 10:50:32.587614",
    "2021-10-05 10:50:32.587617",
    "2021-10-05 10:50:32.587620",
    "2021-10-05 10:50:32.587624",
    "2021-10-05 10:50:32.587627",
    "2021-10-05 10:50:32.587630",
    "2021-10-05 10:50:32.587633",
    "2021-10-05 10:50:32.587636",
    "2021-10-05 10:50:32.587638",
    "2021-10-05 10:50:32.587642",
    "2021-10-05 10:50:32.587645",
    "2021-10-05 10:50:32.587648",
    "2021-10-05 10:50:32.587651",
    "2021-10-05 10:50:32.587654",
    "2021-10-05 10:50:32.587657",
    "2021-10-05 10:50:32.587660",
    "2021-10-05 10:50:32.587663",
    "2021-10-05 10:50:32.587666",
    "2021-10-05 10:50:32.587669"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
